{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per the supplied diff):\n//  1) paragraph \"...avecq le <m>noir d'escaille</m> ou fort,<lb/> ou foible,...\"\n//        -> \"...avecq le <m>noir d'escaille</m>, ou fort,<lb/> ou foible,...\"\n//           (insert a comma immediately after the closing \"</m>\")\n//  2) paragraph \"...de l'<m>esmail d'azur</m> ou fort ou foible, selon qu'ilz veulent<lb/>...\"\n//        -> \"...de l'<m>esmail d'azur</m>, ou fort, ou foible, selon qu'ilz veulent<lb/>...\"\n//           (insert a comma immediately after that closing \"</m>\", AND insert a\n//            comma right after \"ou fort\")\n//\n// Both edits only ever *insert* a single \",\" at a collapsed (zero-length)\n// point, so existing runs/formatting are preserved (a run is split in two,\n// never merged/re-colored).\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Edit 1: \"noir d'escaille</m> ou fort,\"  ->  \"noir d'escaille</m>, ou fort,\"\n// ---------------------------------------------------------------------\n// \"noir d'escaille</m>\" is unique across the whole document, so a body-level\n// search safely lands on the one run that ends in \"</m>\" here.\nconst hit1 = body.search(\"noir d'escaille</m>\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\n\nif (hit1.items.length !== 1) {\n  throw new Error(`Edit 1: expected exactly 1 match for \"noir d'escaille</m>\", found ${hit1.items.length}`);\n}\n// Collapse to the end of the match (right after \"</m>\") and insert \",\" there.\nhit1.items[0].getRange(\"End\").insertText(\",\", \"Before\");\n\n// ---------------------------------------------------------------------\n// Edit 2: \"esmail d'azur</m> ou fort ou foible\" ->\n//         \"esmail d'azur</m>, ou fort, ou foible\"\n// ---------------------------------------------------------------------\n// \"esmail d'azur</m>\" (and \"ou fort\") each appear more than once in the\n// document, so first locate the single paragraph that holds this specific\n// sentence, then scope both searches to that paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"ou fort ou foible\") !== -1) {\n    targetParagraph = p;\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error('Edit 2: could not find the paragraph containing \"ou fort ou foible\"');\n}\n\n// 2a) insert \",\" immediately after \"</m>\"\nconst mHit = targetParagraph.search(\"esmail d'azur</m>\", { matchCase: true });\nmHit.load(\"items\");\nawait context.sync();\nif (mHit.items.length !== 1) {\n  throw new Error(`Edit 2a: expected exactly 1 match for \"esmail d'azur</m>\" in target paragraph, found ${mHit.items.length}`);\n}\nmHit.items[0].getRange(\"End\").insertText(\",\", \"Before\");\n\n// 2b) insert \",\" immediately after \"ou fort\" (re-search so the range\n//     reflects the paragraph after edit 2a's insertion).\nconst fortHit = targetParagraph.search(\"ou fort\", { matchCase: true });\nfortHit.load(\"items\");\nawait context.sync();\nif (fortHit.items.length !== 1) {\n  throw new Error(`Edit 2b: expected exactly 1 match for \"ou fort\" in target paragraph, found ${fortHit.items.length}`);\n}\nfortHit.items[0].getRange(\"End\").insertText(\",\", \"Before\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per the supplied diff):\n#  1) paragraph \"...avecq le <m>noir d'escaille</m> ou fort,<lb/> ou foible,...\"\n#        -> \"...avecq le <m>noir d'escaille</m>, ou fort,<lb/> ou foible,...\"\n#           (insert a comma immediately after the closing \"</m>\")\n#  2) paragraph \"...de l'<m>esmail d'azur</m> ou fort ou foible, selon qu'ilz veulent<lb/>...\"\n#        -> \"...de l'<m>esmail d'azur</m>, ou fort, ou foible, selon qu'ilz veulent<lb/>...\"\n#           (insert a comma immediately after that closing \"</m>\", AND insert a\n#            comma right after \"ou fort\")\n#\n# NOTE: Find.Execute(..., Replace:=...) applies Word's \"smart quotes\"\n# AutoCorrect to the replacement text, which would turn the straight\n# apostrophes in this French text (e.g. \"d'escaille\", \"qu'ilz\") into curly\n# ones -- a change NOT present in the target diff. To avoid that, every edit\n# below only ever uses Find to *locate* a zero-length insertion point and\n# then InsertAfter(\",\") -- a comma has no apostrophe, so AutoCorrect cannot\n# touch anything, and the existing runs/formatting are preserved (a run is\n# only ever split in two, never merged/re-colored).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Edit 1: \"noir d'escaille</m> ou fort,\" -> \"noir d'escaille</m>, ou fort,\"\n# ---------------------------------------------------------------------\n# \"noir d'escaille</m>\" is unique across the whole document, so searching\n# the whole document body safely lands on the one run that ends in \"</m>\"\n# here.\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\"noir d'escaille</m>\")\nif (-not $found1) {\n    throw \"Edit 1: could not find `\"noir d'escaille</m>`\"\"\n}\n$r1.Collapse(0)  # wdCollapseEnd -> collapse to the end of the match\n$r1.InsertAfter(\",\")\n\n# ---------------------------------------------------------------------\n# Edit 2: \"esmail d'azur</m> ou fort ou foible\" ->\n#         \"esmail d'azur</m>, ou fort, ou foible\"\n# ---------------------------------------------------------------------\n# \"esmail d'azur</m>\" (and \"ou fort\") each appear more than once in the\n# document, so first locate the single paragraph that holds this specific\n# sentence, then scope both Find operations to that paragraph's range.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*ou fort ou foible*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Edit 2: could not find the paragraph containing `\"ou fort ou foible`\"\"\n}\n\n# 2a) insert \",\" immediately after \"</m>\"\n$r2a = $target.Range\n$found2a = $r2a.Find.Execute(\"esmail d'azur</m>\")\nif (-not $found2a) {\n    throw \"Edit 2a: could not find `\"esmail d'azur</m>`\" in target paragraph\"\n}\n$r2a.Collapse(0)\n$r2a.InsertAfter(\",\")\n\n# 2b) insert \",\" immediately after \"ou fort\" (re-fetch the paragraph range so\n#     the Find reflects the paragraph's text after edit 2a's insertion).\n$r2b = $target.Range\n$found2b = $r2b.Find.Execute(\"ou fort\")\nif (-not $found2b) {\n    throw \"Edit 2b: could not find `\"ou fort`\" in target paragraph\"\n}\n$r2b.Collapse(0)\n$r2b.InsertAfter(\",\")\n"}
